$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Credentials" -- swap the sample login creds used by the step defs
# ---------------------------------------------------------------------------
$wsCred = $wb.Worksheets.Item("Credentials")
$wsCred.Range("B2").Value = "standard_user"
$wsCred.Range("C2").Value = "secret_sauce"

# ---------------------------------------------------------------------------
# Sheet 2: "GL Reports" -- add a "role" column after scenario_name, reorder
# the line1 description next to the line1 fields, and add two rows of sample
# journal data (one "successful login" + one "unsuccessful login" scenario).
# ---------------------------------------------------------------------------
$wsGl = $wb.Worksheets.Item("GL Reports")

# Insert the new "role" column right after scenario_name; this shifts every
# other header/column (and their bestFit widths) one slot to the right.
$wsGl.Columns("B:B").Insert()
$wsGl.Range("B1").Value = "role"

# The insert leaves journal_line1_desc trailing after journal_line2_credit
# (same relative order as before the insert); put it back next to the other
# line1 fields and shift journal_line2_account/credit/desc over to match.
$wsGl.Range("L1").Value = "journal_line1_desc"
$wsGl.Range("M1").Value = "journal_line2_account"
$wsGl.Range("N1").Value = "journal_line2_credit"
$wsGl.Range("O1").Value = "journal_line2_desc"

# The inserted column has no bestFit width recorded yet -- approximate the
# saved sheet's AutoFit widths for the (now wider) A and new B columns.
$wsGl.Columns("A").ColumnWidth = 34.666666666666664
$wsGl.Columns("B").ColumnWidth = 16.5

# Row 2 - General Accountant / successful login scenario
$wsGl.Range("A2").Value = "Successful login with valid credentials"
$wsGl.Range("B2").Value = "General Accountant"
$wsGl.Range("C2").Value = "abc"
$wsGl.Range("D2").Value = 45925
$wsGl.Range("D2").NumberFormat = "mmm\-dd"
$wsGl.Range("E2").Value = "Spreadsheet"
$wsGl.Range("F2").Value = "PE"
$wsGl.Range("G2").Value = "GBP"
$wsGl.Range("H2").Value = "User"
$wsGl.Range("I2").Value = 1
$wsGl.Range("J2").Value = "1111.1123.1234"
$wsGl.Range("K2").Value = 120
$wsGl.Range("L2").Value = "line1"
$wsGl.Range("M2").Value = "2222.1123.1234"
$wsGl.Range("N2").Value = 120
$wsGl.Range("O2").Value = "line2"

# Row 3 - Approver / unsuccessful login scenario
$wsGl.Range("A3").Value = "Unsuccessful login with invalid credentials"
$wsGl.Range("B3").Value = "Approver"
$wsGl.Range("C3").Value = "abc"
$wsGl.Range("D3").Value = 45925
$wsGl.Range("D3").NumberFormat = "mmm\-dd"
$wsGl.Range("E3").Value = "Spreadsheet"
$wsGl.Range("F3").Value = "PE"
$wsGl.Range("G3").Value = "GBP"
$wsGl.Range("H3").Value = "User"
$wsGl.Range("I3").Value = 1
$wsGl.Range("J3").Value = "1111.1123.1234"
$wsGl.Range("K3").Value = 120
$wsGl.Range("L3").Value = "line1"
$wsGl.Range("M3").Value = "2222.1123.1234"
$wsGl.Range("N3").Value = 120
$wsGl.Range("O3").Value = "line2"

# Row 4 - a leftover formatted-but-empty date cell under the date column
$wsGl.Range("D4").NumberFormat = "mmm\-dd"

# ---------------------------------------------------------------------------
# Selection / active sheet bookkeeping to mirror the saved view state
# ---------------------------------------------------------------------------
$wsGl.Range("B7").Select()
$wsCred.Activate()
$wsCred.Range("C6").Select()
